# Apply updated cryptocurrency price/volume figures to Sheet1.
# D-column price values are forced to Text via a leading quote-prefix
# (Excel's standard 'treat as text' convention) so that values such as
# "535.80" or "58.751.23" keep their exact literal formatting instead of
# being auto-coerced into numbers (which would drop trailing zeros / fail
# to parse dot-grouped thousands). E-column percentages already contain
# surrounding spaces, so Excel stores them as text without any extra help.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.751.23"
$ws.Range("E2").Value = "  -1.06%  "
$ws.Range("D3").Value = "'2.518.82"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'535.80"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("D6").Value = "'136.55"
$ws.Range("E6").Value = "  -2.14%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("E8").Value = "  +0.54%  "
$ws.Range("D9").Value = "'2.517.17"
$ws.Range("E9").Value = "  -0.16%  "
$ws.Range("E10").Value = "  -0.65%  "
$ws.Range("E11").Value = "  -2.46%  "
$ws.Range("E12").Value = "  -1.65%  "
$ws.Range("E13").Value = "  -3.52%  "
$ws.Range("D14").Value = "'2.935.40"
$ws.Range("E14").Value = "  -1.02%  "
$ws.Range("D15").Value = "'22.86"
$ws.Range("E15").Value = "  -2.61%  "
$ws.Range("D16").Value = "'58.676.04"
$ws.Range("E16").Value = "  -1.02%  "
$ws.Range("E17").Value = "  -1.76%  "
$ws.Range("D18").Value = "'2.508.49"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("D19").Value = "'11.06"
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("D20").Value = "'4.26"
$ws.Range("E20").Value = "  -1.24%  "
$ws.Range("D21").Value = "'322.87"
$ws.Range("E21").Value = "  -0.93%  "
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").Value = "'5.91"
$ws.Range("E23").Value = "  +1.29%  "
$ws.Range("D24").Value = "'65.35"
$ws.Range("E24").Value = "  +3.37%  "
$ws.Range("D25").Value = "'0.419"
$ws.Range("E25").Value = "  -1.42%  "
$ws.Range("E26").Value = "  -1.54%  "
$ws.Range("D27").Value = "'0.997"
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("D28").Value = "'7.54"
$ws.Range("E28").Value = "  -3.66%  "
$ws.Range("D29").Value = "'6.68"
$ws.Range("E29").Value = "  -3.94%  "
$ws.Range("D30").Value = "'0.0₃0765"
$ws.Range("E30").Value = "  -1.94%  "
$ws.Range("E31").Value = "  -1.72%  "
$ws.Range("D32").Value = "'166.60"
$ws.Range("E32").Value = "  +1.02%  "
$ws.Range("E33").Value = "  +4.31%  "
$ws.Range("D34").Value = "'0.998"
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("D35").Value = "'1.45"
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("D36").Value = "'18.43"
$ws.Range("E36").Value = "  -0.55%  "
$ws.Range("E37").Value = "  -4.33%  "
$ws.Range("E38").Value = "  -3.56%  "
$ws.Range("E39").Value = "  -0.87%  "
$ws.Range("D40").Value = "'0.812"
$ws.Range("E40").Value = "  -0.12%  "
$ws.Range("D41").Value = "'3.60"
$ws.Range("E41").Value = "  -2.36%  "
$ws.Range("D42").Value = "'284.19"
$ws.Range("E42").Value = "  +1.59%  "
$ws.Range("E43").Value = "  -0.72%  "
$ws.Range("D44").Value = "'0.997"
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").Value = "'130.89"
$ws.Range("E45").Value = "  +5.95%  "
$ws.Range("D46").Value = "'0.604"
$ws.Range("E46").Value = "  +0.70%  "
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("E48").Value = "  -1.04%  "
$ws.Range("E49").Value = "  -1.96%  "
$ws.Range("E50").Value = "  -2.36%  "
$ws.Range("D51").Value = "'17.18"
$ws.Range("E51").Value = "  -3.87%  "
